# Generate Report for handback
# Updates the Overview status text and records the "handed back" details
# (target file / handback file / handback datetime / handoff reason) for
# the zh-cn and de-de localization reports.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: both language columns move from "Ready for handoff"
# to "Handed back: in sync with en-US" for the two tracked files.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("B2").Value = $newStatus
$wsZh.Range("B3").Value = $newStatus

# Row 2 (05daeeb2-...md)
$wsZh.Range("E2").Value = "05daeeb2-ac12-4594-84bc-3e5a63870673.md"
$wsZh.Hyperlinks.Add($wsZh.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/af756787332082ec5e6400dc07cae1b2cb88645d/e2e/05daeeb2-ac12-4594-84bc-3e5a63870673.md", "", "", "05daeeb2-ac12-4594-84bc-3e5a63870673.md") | Out-Null

$wsZh.Range("F2").Value = "05daeeb2-ac12-4594-84bc-3e5a63870673.25c69d1c765d57f5d8d6eccac1dd13afc922757f.zh-cn.xlf"
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/745a7c469fbbcbc0aa50969d48f5ba3e387b4693/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/05daeeb2-ac12-4594-84bc-3e5a63870673.25c69d1c765d57f5d8d6eccac1dd13afc922757f.zh-cn.xlf", "", "", "05daeeb2-ac12-4594-84bc-3e5a63870673.25c69d1c765d57f5d8d6eccac1dd13afc922757f.zh-cn.xlf") | Out-Null

$wsZh.Range("G2").Value = "2016-01-14 03:13:59"

# Row 3 (097298bd-...md)
$wsZh.Range("E3").Value = "097298bd-7ea0-4fd0-9d70-83728cc19d14.md"
$wsZh.Hyperlinks.Add($wsZh.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/af756787332082ec5e6400dc07cae1b2cb88645d/e2e/097298bd-7ea0-4fd0-9d70-83728cc19d14.md", "", "", "097298bd-7ea0-4fd0-9d70-83728cc19d14.md") | Out-Null

$wsZh.Range("F3").Value = "097298bd-7ea0-4fd0-9d70-83728cc19d14.f32b93eeefa6fc8adad04f453fa5c24b491312ad.zh-cn.xlf"
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/745a7c469fbbcbc0aa50969d48f5ba3e387b4693/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/097298bd-7ea0-4fd0-9d70-83728cc19d14.f32b93eeefa6fc8adad04f453fa5c24b491312ad.zh-cn.xlf", "", "", "097298bd-7ea0-4fd0-9d70-83728cc19d14.f32b93eeefa6fc8adad04f453fa5c24b491312ad.zh-cn.xlf") | Out-Null

$wsZh.Range("G3").Value = "2016-01-14 03:13:59"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("B2").Value = $newStatus
$wsDe.Range("B3").Value = $newStatus

# Row 2 (05daeeb2-...md)
$wsDe.Range("E2").Value = "05daeeb2-ac12-4594-84bc-3e5a63870673.md"
$wsDe.Hyperlinks.Add($wsDe.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/af756787332082ec5e6400dc07cae1b2cb88645d/e2e/05daeeb2-ac12-4594-84bc-3e5a63870673.md", "", "", "05daeeb2-ac12-4594-84bc-3e5a63870673.md") | Out-Null

$wsDe.Range("F2").Value = "05daeeb2-ac12-4594-84bc-3e5a63870673.25c69d1c765d57f5d8d6eccac1dd13afc922757f.de-de.xlf"
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4675af6c2255ee046c2bb2d20562c04aa51a6779/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/05daeeb2-ac12-4594-84bc-3e5a63870673.25c69d1c765d57f5d8d6eccac1dd13afc922757f.de-de.xlf", "", "", "05daeeb2-ac12-4594-84bc-3e5a63870673.25c69d1c765d57f5d8d6eccac1dd13afc922757f.de-de.xlf") | Out-Null

$wsDe.Range("G2").Value = "2016-01-14 03:14:21"

# Row 3 (097298bd-...md)
$wsDe.Range("E3").Value = "097298bd-7ea0-4fd0-9d70-83728cc19d14.md"
$wsDe.Hyperlinks.Add($wsDe.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/af756787332082ec5e6400dc07cae1b2cb88645d/e2e/097298bd-7ea0-4fd0-9d70-83728cc19d14.md", "", "", "097298bd-7ea0-4fd0-9d70-83728cc19d14.md") | Out-Null

$wsDe.Range("F3").Value = "097298bd-7ea0-4fd0-9d70-83728cc19d14.f32b93eeefa6fc8adad04f453fa5c24b491312ad.de-de.xlf"
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4675af6c2255ee046c2bb2d20562c04aa51a6779/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/097298bd-7ea0-4fd0-9d70-83728cc19d14.f32b93eeefa6fc8adad04f453fa5c24b491312ad.de-de.xlf", "", "", "097298bd-7ea0-4fd0-9d70-83728cc19d14.f32b93eeefa6fc8adad04f453fa5c24b491312ad.de-de.xlf") | Out-Null

$wsDe.Range("G3").Value = "2016-01-14 03:14:21"

Write-Host "Handback report generated."
